# Apply updated odds values for 2025-12-31 Betfair Back/Lay sheet.
# Maps directly onto the cell-level diff between the prior and current
# export of the workbook (single worksheet, row-per-fixture layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 3.45
$ws.Range("AC2").Value = 8
$ws.Range("AK2").Value = 1000
$ws.Range("R5").Value = 1.27
$ws.Range("F6").Value = 6.4
$ws.Range("G6").Value = 9.199999999999999
$ws.Range("H6").Value = 1.39
$ws.Range("I6").Value = 1.53
$ws.Range("J6").Value = 4.9
$ws.Range("K6").Value = 7.4
$ws.Range("N6").Value = 4.4
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 2.3
$ws.Range("Q6").Value = 1.58
$ws.Range("R6").Value = 1.51
$ws.Range("S6").Value = 2.5
$ws.Range("T6").Value = 1.81
$ws.Range("U6").Value = 1.94
$ws.Range("V6").Value = 2.88
$ws.Range("W6").Value = 1.13
$ws.Range("AF6").Value = 80
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 2.54
$ws.Range("Q7").Value = 1.42
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 2.04
$ws.Range("V7").Value = 1.78
$ws.Range("X7").Value = 34
$ws.Range("Y7").Value = 16.5
$ws.Range("Z7").Value = 18.5
$ws.Range("AA7").Value = 29
$ws.Range("AB7").Value = 22
$ws.Range("AC7").Value = 11
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 30
$ws.Range("AG7").Value = 16
$ws.Range("AK7").Value = 34
$ws.Range("AL7").Value = 36
$ws.Range("AM7").Value = 55
$ws.Range("AN7").Value = 20
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 3.4
$ws.Range("L8").Value = 1.24
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 5.1
$ws.Range("O8").Value = 1.19
$ws.Range("R8").Value = 1.56
$ws.Range("S8").Value = 2.4
$ws.Range("T8").Value = 1.54
$ws.Range("U8").Value = 2.48
$ws.Range("V8").Value = 1.39
$ws.Range("W8").Value = 1.73
$ws.Range("X8").Value = 29
$ws.Range("Y8").Value = 22
$ws.Range("Z8").Value = 28
$ws.Range("AA8").Value = 65
$ws.Range("AB8").Value = 17.5
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 15.5
$ws.Range("AE8").Value = 40
$ws.Range("AF8").Value = 21
$ws.Range("AG8").Value = 14
$ws.Range("AH8").Value = 18.5
$ws.Range("AI8").Value = 38
$ws.Range("AJ8").Value = 36
$ws.Range("AK8").Value = 26
$ws.Range("AL8").Value = 980
$ws.Range("AM8").Value = 70
$ws.Range("AN8").Value = 14
$ws.Range("AO8").Value = 26
$ws.Range("F9").Value = 5.3
$ws.Range("G9").Value = 6.4
$ws.Range("H9").Value = 1.64
$ws.Range("I9").Value = 1.78
$ws.Range("K9").Value = 4.5
$ws.Range("Q9").Value = 1.78
$ws.Range("U9").Value = 1.98
$ws.Range("V9").Value = 2.28
$ws.Range("W9").Value = 1.19
$ws.Range("X9").Value = 990
$ws.Range("AC9").Value = 10
$ws.Range("T10").Value = 1.55
$ws.Range("F11").Value = 4.8
$ws.Range("G11").Value = 6.8
$ws.Range("H11").Value = 1.62
$ws.Range("I11").Value = 1.78
$ws.Range("J11").Value = 4.1
$ws.Range("K11").Value = 4.9
$ws.Range("N11").Value = 2.3
$ws.Range("P11").Value = 2.3
$ws.Range("Q11").Value = 1.6
$ws.Range("V11").Value = 2.28
$ws.Range("W11").Value = 1.17
$ws.Range("F12").Value = 1.27
$ws.Range("G12").Value = 1.35
$ws.Range("H12").Value = 9.4
$ws.Range("I12").Value = 16
$ws.Range("J12").Value = 5.9
$ws.Range("K12").Value = 8.4
$ws.Range("L12").Value = 1.18
$ws.Range("N12").Value = 5.6
$ws.Range("O12").Value = 1.16
$ws.Range("P12").Value = 2.6
$ws.Range("Q12").Value = 1.48
$ws.Range("R12").Value = 1.64
$ws.Range("S12").Value = 2.2
$ws.Range("T12").Value = 1.92
$ws.Range("U12").Value = 1.85
$ws.Range("V12").Value = 1.06
$ws.Range("W12").Value = 3.8
$ws.Range("AB12").Value = 16
$ws.Range("AF12").Value = 11
$ws.Range("AG12").Value = 13.5
$ws.Range("AJ12").Value = 13
$ws.Range("AK12").Value = 17
